$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TPSEE")
$ws.Activate()
$ws.Range("C2").Value = "US"
$ws.Range("D2").Select()
